$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 268; this shifts the existing
# rows 268-276 down to 269-277 and extends the used range to R277.
$ws.Rows.Item(268).Insert()

# Populate the newly inserted row 268 with the new weekly record.
$ws.Range("A268").Value = 4
$ws.Range("B268").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C268").Value = "Los Lagos"
$ws.Range("D268").Value = 44747
$ws.Range("E268").Value = 10
$ws.Range("F268").Value = 100112021
$ws.Range("G268").Value = "Ají"
$ws.Range("H268").Value = "Inferno"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 160
$ws.Range("K268").Value = 18000
$ws.Range("L268").Value = 19000
$ws.Range("M268").Value = 18500
$ws.Range("N268").Value = "$/caja 12 kilos"
$ws.Range("O268").Value = "Región de Arica y Parinacota"
$ws.Range("P268").Value = 1542
$ws.Range("Q268").Value = 12
$ws.Range("R268").Value = "Hortaliza"
